$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-27 23:18:18"
$ws.Range("H2").Value = "'57%"
$ws.Range("E3").Value = "2026-02-27 23:18:21"
$ws.Range("E4").Value = "2026-02-27 23:18:23"
$ws.Range("E5").Value = "2026-02-27 23:18:26"
$ws.Range("H5").Value = "'48%"
$ws.Range("O5").Value = "4.4 °C"
$ws.Range("E6").Value = "2026-02-27 23:18:28"
$ws.Range("E7").Value = "2026-02-27 23:18:30"
$ws.Range("E8").Value = "2026-02-27 23:18:33"
$ws.Range("H8").Value = "'69%"
$ws.Range("N8").Value = "8.0 °C 22:30 TU"
$ws.Range("O8").Value = "11.4 °C"
$ws.Range("E9").Value = "2026-02-27 23:18:36"
$ws.Range("K9").Value = "7.9 MJ/m2"
$ws.Range("O9").Value = "11.0 °C"
$ws.Range("E10").Value = "2026-02-27 23:18:38"
$ws.Range("H10").Value = "'87%"
$ws.Range("E11").Value = "2026-02-27 23:18:41"
$ws.Range("H11").Value = "'73%"
$ws.Range("O11").Value = "8.3 °C"
$ws.Range("E12").Value = "2026-02-27 23:18:43"
$ws.Range("H12").Value = "'96%"
$ws.Range("O12").Value = "10.7 °C"
$ws.Range("E13").Value = "2026-02-27 23:18:45"
$ws.Range("H13").Value = "'66%"
$ws.Range("O13").Value = "6.5 °C"
$ws.Range("E14").Value = "2026-02-27 23:18:48"
$ws.Range("N14").Value = "5.3 °C 22:24 TU"
$ws.Range("O14").Value = "10.2 °C"
$ws.Range("E15").Value = "2026-02-27 23:18:50"
$ws.Range("E16").Value = "2026-02-27 23:18:53"
$ws.Range("N16").Value = "-0.6 °C 22:43 TU"
$ws.Range("E17").Value = "2026-02-27 23:18:55"
$ws.Range("O17").Value = "7.3 °C"
$ws.Range("E18").Value = "2026-02-27 23:18:58"
$ws.Range("N18").Value = "8.9 °C 22:59 TU"
$ws.Range("O18").Value = "11.8 °C"
$ws.Range("E19").Value = "2026-02-27 23:19:00"
$ws.Range("N19").Value = "6.7 °C 22:58 TU"
$ws.Range("E20").Value = "2026-02-27 23:19:02"
$ws.Range("E21").Value = "2026-02-27 23:19:05"
$ws.Range("O21").Value = "9.7 °C"
$ws.Range("E22").Value = "2026-02-27 23:19:07"
$ws.Range("E23").Value = "2026-02-27 23:19:10"
$ws.Range("H23").Value = "'44%"
$ws.Range("N23").Value = "0.7 °C 22:44 TU"
$ws.Range("E24").Value = "2026-02-27 23:19:12"
$ws.Range("O24").Value = "10.0 °C"
$ws.Range("E25").Value = "2026-02-27 23:19:15"
$ws.Range("O25").Value = "5.7 °C"
$ws.Range("E26").Value = "2026-02-27 23:19:17"
$ws.Range("H26").Value = "'48%"
$ws.Range("J26").Value = "1021.6 hPa"
$ws.Range("N26").Value = "4.9 °C 22:28 TU"
$ws.Range("O26").Value = "9.8 °C"
$ws.Range("E27").Value = "2026-02-27 23:19:19"
$ws.Range("E28").Value = "2026-02-27 23:19:22"
$ws.Range("E29").Value = "2026-02-27 23:19:24"
$ws.Range("O29").Value = "11.6 °C"
$ws.Range("E30").Value = "2026-02-27 23:19:27"
$ws.Range("H30").Value = "'92%"
$ws.Range("N30").Value = "8.6 °C 22:59 TU"
$ws.Range("O30").Value = "10.8 °C"
$ws.Range("E31").Value = "2026-02-27 23:19:29"
$ws.Range("E32").Value = "2026-02-27 23:19:31"
$ws.Range("H32").Value = "'54%"
$ws.Range("E33").Value = "2026-02-27 23:19:34"
$ws.Range("J33").Value = "1023.6 hPa"
$ws.Range("O33").Value = "8.5 °C"
$ws.Range("E34").Value = "2026-02-27 23:19:36"
$ws.Range("H34").Value = "'49%"
$ws.Range("O34").Value = "4.4 °C"
$ws.Range("E35").Value = "2026-02-27 23:19:39"
$ws.Range("N35").Value = "7.5 °C 22:47 TU"
$ws.Range("O35").Value = "11.8 °C"
$ws.Range("E36").Value = "2026-02-27 23:19:41"
$ws.Range("E37").Value = "2026-02-27 23:19:44"
$ws.Range("H37").Value = "'70%"
$ws.Range("E38").Value = "2026-02-27 23:19:46"
$ws.Range("E39").Value = "2026-02-27 23:19:49"
$ws.Range("E40").Value = "2026-02-27 23:19:52"
$ws.Range("O40").Value = "8.8 °C"
$ws.Range("E41").Value = "2026-02-27 23:19:54"
$ws.Range("E42").Value = "2026-02-27 23:19:57"
$ws.Range("O42").Value = "11.2 °C"
$ws.Range("E43").Value = "2026-02-27 23:19:59"
$ws.Range("E44").Value = "2026-02-27 23:20:02"
$ws.Range("H44").Value = "'62%"
$ws.Range("E45").Value = "2026-02-27 23:20:04"
$ws.Range("H45").Value = "'47%"
$ws.Range("O45").Value = "11.6 °C"
$ws.Range("E46").Value = "2026-02-27 23:20:07"
$ws.Range("J46").Value = "1023.9 hPa"
